$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-style row 43: it used to be the last ("isolated") entry in the table
#    (top+bottom thin border, style indices 8/9). Now that two more rows are
#    appended after it, row 43 becomes a normal/interior row, and the new
#    row 40 ... no wait: row 40 already carries the "isolated" look, which is
#    exactly what row 43 should look like after the edit, so copy its format.
# ---------------------------------------------------------------------------
$ws.Range("A40:E40").Copy()
$ws.Range("A43:E43").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Prepare formatting for the two new rows by cloning it from existing rows
#    that already carry the right look:
#      - row 44 (has a value in column A) <- row 41's format
#      - row 45 (no value in column A)    <- row 9's format (columns B:E)
# ---------------------------------------------------------------------------
$ws.Range("A41:E41").Copy()
$ws.Range("A44:E44").PasteSpecial(-4122)

$ws.Range("B9:E9").Copy()
$ws.Range("B45:E45").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill in the new cell values. The order below matches the order the new
#    strings were first introduced into the shared-string table.
# ---------------------------------------------------------------------------
$ws.Range("C44").Value = " Wasn\'t the expedition fun?!"
$ws.Range("C45").Value = " Well, we\'re back to the regular\nguild routine! It\'s time to work again!"
$ws.Range("A44").Value = "SCRIPT/G01P03A/um1101.ssb"
$ws.Range("D44").Value = " Разве в нашей экспедиции не\nбыло весело?!"
$ws.Range("D45").Value = " Что-ж, теперь мы снова\nзанимаемся нашими гильдейскими делами!\nПора снова поработать!"
$ws.Range("E44").Value = " Ñàèâå â îàšåê üëòðåäéøéé îå\náúìï âåòåìï?!"
$ws.Range("E45").Value = " Œóï-ç, óåðåñû íú òîïâà\nèàîéíàåíòÿ îàšéíé ãéìûäåêòëéíé äåìàíé!\nÐïñà òîïâà ðïñàáïóàóû!"

$ws.Range("B44").Value = 570
$ws.Range("B45").Value = 573

# ---------------------------------------------------------------------------
# 4. Row heights for the two new rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(44).RowHeight = 43.2
$ws.Rows.Item(45).RowHeight = 32.4

# ---------------------------------------------------------------------------
# 5. Scroll / selection bookkeeping to match where the editor ended up.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D45").Select()
